{"js": "const body = context.document.body;\n\n// --- Edit 1: merge the two \"THU Nov 22\" / \" 14:59:44 IST 2018\" runs into one run ---\nconst dateSearch = body.search(\"THU Nov 22\", { matchCase: true });\ndateSearch.load(\"items\");\nawait context.sync();\n\nif (dateSearch.items.length === 0) {\n  throw new Error(\"Could not find 'THU Nov 22' text to merge runs.\");\n}\n\nconst dateParagraph = dateSearch.items[0].paragraphs.getFirst();\nconst dateRange = dateParagraph.getRange();\ndateRange.insertText(\"THU Nov 22 14:59:44 IST 2018\", Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Edit 2: append a new \"chick in\" purchase-details block after the existing\n// final \"Amount balance\" entry (the one that reads \"- 920.0\") ---\nconst balanceSearch = body.search(\"920.0\", { matchCase: true });\nbalanceSearch.load(\"items\");\nawait context.sync();\n\nif (balanceSearch.items.length < 2) {\n  throw new Error(\"Could not find the second '920.0' occurrence (Amount balance row).\");\n}\n\nconst anchorParagraph = balanceSearch.items[balanceSearch.items.length - 1].paragraphs.getFirst();\nconst anchorRange = anchorParagraph.getRange();\n\nconst newBlockOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n<w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n      </w:pPr>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>FRI Nov 23</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t xml:space=\"preserve\"> 11:13:49 IST 2018</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Person Name</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- HAUMANTHARAYA</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Bill number</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- 9065</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>---------------------------------------------------------------</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Item Name</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- NAVILKOSU</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Number of Pockets</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- 1</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Number of KGs</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- 63</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Rate</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- 20</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:t>Total Price</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- 1260.0</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n      </w:pPr>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n        <w:t>Amount balance</w:t>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n        <w:tab/>\n      </w:r>\n      <w:r>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n        <w:tab/>\n        <w:t>- 1260.0</w:t>\n      </w:r>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n        </w:rPr>\n      </w:pPr>\n    </w:p><w:p>\n      <w:pPr>\n        <w:pStyle w:val=\"PlainText\"/>\n        <w:rPr>\n          <w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/>\n          <w:b/>\n        </w:rPr>\n      </w:pPr>\n     </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nanchorRange.insertOoxml(newBlockOoxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1: merge the two \"THU Nov 22\" / \" 14:59:44 IST 2018\" runs into one run ---\n$mergedDate = \"THU Nov 22 14:59:44 IST 2018\"\n$found = $d.Content.Find.Execute($mergedDate, $false, $false, $false, $false, $false, $true, 1, $false, $mergedDate, 2)\nif (-not $found) {\n    throw \"Could not find the THU Nov 22 date paragraph to merge.\"\n}\n\n# --- Edit 2: append a new chick-in purchase-details block after the existing\n# final Amount balance entry (the one that reads \"- 920.0\") ---\n$anchorPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*920.0*\") {\n        $anchorPara = $p\n    }\n}\nif (-not $anchorPara) {\n    throw \"Could not find the Amount balance paragraph to anchor the new block.\"\n}\n\n$insertRng = $anchorPara.Range\n$insertRng.Collapse(0)\n\n$newBlockXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>FRI Nov 23</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t xml:space=\"preserve\"> 11:13:49 IST 2018</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Person Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- HAUMANTHARAYA</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Bill number</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- 9065</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>---------------------------------------------------------------</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Item Name</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- NAVILKOSU</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Number of Pockets</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- 1</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Number of KGs</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- 63</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Rate</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- 20</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:t>Total Price</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr><w:tab/><w:t>- 1260.0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr><w:t>Amount balance</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr><w:tab/><w:t>- 1260.0</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val=\"PlainText\"/><w:rPr><w:rFonts w:ascii=\"Courier New\" w:hAnsi=\"Courier New\" w:cs=\"Courier New\"/><w:b/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n$insertRng.InsertXML($newBlockXml)\n"}
